# Update countries & provincias Spain
# Refresh the COVID-19 "paises" table with the latest report figures.
#
# The underlying data source re-sorted a few countries that share the same
# "Casos totales" ranking (Barein/Bolivia/Israel, Islas Malvinas/Groenlandia,
# Santa Sede/Islas Turcas y Caicos), which is why some rows below change the
# country name in column A as well as their statistics. Most rows simply get
# refreshed numbers for the same country.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 17 de Junio de 2020 a las 23:57"

# --- Pure numeric refreshes (country in column A is unchanged) -----------

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 2231581
$ws.Cells.Item(4, 3).Value = 23181
$ws.Cells.Item(4, 4).Value = 911110
$ws.Cells.Item(4, 5).Value = 1200580
$ws.Cells.Item(4, 7).Value = 759
$ws.Cells.Item(4, 8).Value = 119891

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 955377
$ws.Cells.Item(5, 3).Value = 26543
$ws.Cells.Item(5, 5).Value = 431503
$ws.Cells.Item(5, 7).Value = 1054
$ws.Cells.Item(5, 8).Value = 46510

# Row 13: Alemania
$ws.Cells.Item(13, 2).Value = 190179
$ws.Cells.Item(13, 3).Value = 1797
$ws.Cells.Item(13, 5).Value = 7652
$ws.Cells.Item(13, 7).Value = 17
$ws.Cells.Item(13, 8).Value = 8927

# Row 29: Egipto
$ws.Cells.Item(29, 4).Value = 13141
$ws.Cells.Item(29, 5).Value = 34228

# Row 90: Bulgaria
$ws.Cells.Item(90, 2).Value = 3542
$ws.Cells.Item(90, 3).Value = 89
$ws.Cells.Item(90, 4).Value = 1880
$ws.Cells.Item(90, 5).Value = 1478
$ws.Cells.Item(90, 7).Value = 3
$ws.Cells.Item(90, 8).Value = 184

# Row 148: Estado de Palestina
$ws.Cells.Item(148, 2).Value = 555
$ws.Cells.Item(148, 3).Value = 41
$ws.Cells.Item(148, 5).Value = 137

# Row 159: Montenegro
$ws.Cells.Item(159, 2).Value = 333
$ws.Cells.Item(159, 3).Value = 7
$ws.Cells.Item(159, 5).Value = 9

# Row 161: Surinam
$ws.Cells.Item(161, 2).Value = 242
$ws.Cells.Item(161, 3).Value = 6
$ws.Cells.Item(161, 5).Value = 188

# Row 181: Barbados
$ws.Cells.Item(181, 4).Value = 85
$ws.Cells.Item(181, 5).Value = 5

# Row 183: Botsuana
$ws.Cells.Item(183, 4).Value = 25
$ws.Cells.Item(183, 5).Value = 53

# --- Rows whose country re-ranked, so column A also changes ---------------

# Row 49 becomes Barein (was Bolivia)
$ws.Cells.Item(49, 1).Value = "Barein"
$ws.Cells.Item(49, 2).Value = 19961
$ws.Cells.Item(49, 3).Value = 408
$ws.Cells.Item(49, 4).Value = 14185
$ws.Cells.Item(49, 5).Value = 5727
$ws.Cells.Item(49, 7).Value = 2
$ws.Cells.Item(49, 8).Value = 49

# Row 50 becomes Bolivia (was Israel)
$ws.Cells.Item(50, 1).Value = "Bolivia"
$ws.Cells.Item(50, 2).Value = 19883
$ws.Cells.Item(50, 3).Value = 810
$ws.Cells.Item(50, 4).Value = 3752
$ws.Cells.Item(50, 5).Value = 15472
$ws.Cells.Item(50, 7).Value = 27
$ws.Cells.Item(50, 8).Value = 659

# Row 51 becomes Israel (was Barein)
$ws.Cells.Item(51, 1).Value = "Israel"
$ws.Cells.Item(51, 2).Value = 19783
$ws.Cells.Item(51, 3).Value = 288
$ws.Cells.Item(51, 4).Value = 15459
$ws.Cells.Item(51, 5).Value = 4021
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 303

# Row 206 becomes Islas Malvinas (was Groenlandia); stats unchanged
$ws.Cells.Item(206, 1).Value = "Islas Malvinas"

# Row 207 becomes Groenlandia (was Islas Malvinas); stats unchanged
$ws.Cells.Item(207, 1).Value = "Groenlandia"

# Row 208 becomes Santa Sede (was Islas Turcas y Caicos)
$ws.Cells.Item(208, 1).Value = "Santa Sede"
$ws.Cells.Item(208, 4).Value = 12
$ws.Cells.Item(208, 8).Value = 0

# Row 209 becomes Islas Turcas y Caicos (was Santa Sede)
$ws.Cells.Item(209, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 8).Value = 1
